$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")
$ws.Activate()

# --- Row 4 (Question_3 / Amazon Bloomberg screenshot question) ---
# C4: question text is trimmed - the trailing image-reference line
#     ("' images/Options/LongCall.JPG") is removed, a trailing space remains.
$ws.Range("C4").Value = "D'après ce screenshot Bloomberg, que devrait-être le delta de ce call sur une action Amazon Inc ne payant pas de dividendes : "

# D4: the placeholder label changes from "ImageinQuestion" to "ImageinQuestion_pricing"
$ws.Range("D4").Value = "ImageinQuestion_pricing"

# G4 keeps its existing value ("images/Bloomberg/Amazon_call.png") - unchanged.

# --- New column H: "Parameters" ---
# H1 header, formatted like the rest of the header row (style copied from G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Parameters"

# H4 holds the pricing parameters for the new question, formatted like D4/F4 (style 20)
$ws.Range("F4").Copy()
$ws.Range("H4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H4").Value = "132.21, 132.21, 92/360, 0.33532, 0.05363"

$excel.CutCopyMode = 0

# --- View state: scroll / selection moved ---
$ws.Range("J2").Select()
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1

Write-Host "done"
